# Apply updated "想去人数" (F column) values per sheet, per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 15
$ws.Range("F5").Value = 1017
$ws.Range("F7").Value = 2568
$ws.Range("F9").Value = 1257
$ws.Range("F10").Value = 910
$ws.Range("F11").Value = 608
$ws.Range("F12").Value = 916
$ws.Range("F13").Value = 1145
$ws.Range("F16").Value = 114
$ws.Range("F17").Value = 733
$ws.Range("F18").Value = 781
$ws.Range("F19").Value = 199
$ws.Range("F20").Value = 499
$ws.Range("F21").Value = 1121
$ws.Range("F22").Value = 93
$ws.Range("F23").Value = 613
$ws.Range("F24").Value = 596
$ws.Range("F26").Value = 302
$ws.Range("F27").Value = 302
$ws.Range("F29").Value = 475
$ws.Range("F30").Value = 4476
$ws.Range("F31").Value = 487
$ws.Range("F33").Value = 300
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 1613
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 42
$ws.Range("F39").Value = 440
$ws.Range("F42").Value = 143
$ws.Range("F43").Value = 71
$ws.Range("F45").Value = 131
$ws.Range("F47").Value = 110

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 21
$ws.Range("F10").Value = 179
$ws.Range("F12").Value = 189
$ws.Range("F14").Value = 22
$ws.Range("F16").Value = 29
$ws.Range("F17").Value = 191
$ws.Range("F22").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2271

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2271
$ws.Range("F4").Value = 15
$ws.Range("F7").Value = 1017
$ws.Range("F8").Value = 2568
$ws.Range("F10").Value = 1257
$ws.Range("F11").Value = 910
$ws.Range("F12").Value = 608
$ws.Range("F13").Value = 916
$ws.Range("F14").Value = 1145
$ws.Range("F17").Value = 114
$ws.Range("F18").Value = 733
$ws.Range("F19").Value = 781
$ws.Range("F20").Value = 199
$ws.Range("F21").Value = 499
$ws.Range("F22").Value = 1121
$ws.Range("F24").Value = 93
$ws.Range("F25").Value = 613
$ws.Range("F26").Value = 596
$ws.Range("F28").Value = 302
$ws.Range("F30").Value = 476
$ws.Range("F31").Value = 4476
$ws.Range("F32").Value = 189
$ws.Range("F33").Value = 487
$ws.Range("F36").Value = 154
$ws.Range("F37").Value = 1613
$ws.Range("F38").Value = 440
$ws.Range("F39").Value = 23
$ws.Range("F40").Value = 23
$ws.Range("F43").Value = 143
$ws.Range("F44").Value = 71
$ws.Range("F46").Value = 131
$ws.Range("F48").Value = 110
